$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual cell values (shared string text changes)
$ws.Range("M2").Value = "IX"
$ws.Range("V4").Value = "K"

# Clear contents (but keep formatting/style) of rows 5 and 6
$ws.Range("A5:V6").ClearContents()

# Update the view: scrolled to column J, active cell V5 selected
$ws.Range("V5").Select()
$excel.ActiveWindow.ScrollColumn = 10
